$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the header row (row 1: "DIV 1", "DIV 2", "DIV 3"), shifting all
# player data up by one row.
$ws.Rows.Item(1).Delete()

# Reset the selection to match the state left behind by deleting a whole
# row in Excel (the newly-promoted row 1 becomes selected in full).
$ws.Rows.Item(1).Select() | Out-Null
